$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iteration-6")
$ws.Range("A1").Value2 = "STORY"
$ws.Range("A2").Value2 = "STORY"
